$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Re-write existing numeric cells so they serialize with full double
#     precision (matches a genuine Excel re-save / full recalc of the
#     original values - the underlying numbers are unchanged, only their
#     textual representation gains precision). ---
$ws.Range("Q2").Value2 = 1.0793809267402219
$ws.Range("R2").Value2 = 30.221667564092229
$ws.Range("S2").Value2 = 2.9428571428571431
$ws.Range("T2").Value2 = 2.9428571428571431
$ws.Range("U2").Value2 = 1.0793809267402219
$ws.Range("V2").Value2 = 30.221667564092229

$ws.Range("N3").Value2 = 0.96470588235294119
$ws.Range("Q3").Value2 = 1.3437347467010949
$ws.Range("R3").Value2 = 44.687957598967159
$ws.Range("T3").Value2 = 5.0999999999999996
$ws.Range("V3").Value2 = 49.415189205394398
$ws.Range("Z3").Value2 = 1.0034482758620691

$ws.Range("N4").Value2 = 0.98969072164948457
$ws.Range("Q4").Value2 = 1.3488367985849901
$ws.Range("S4").Value2 = 3.8529411764705879
$ws.Range("U4").Value2 = 1.4101998819734449
$ws.Range("Z4").Value2 = 1.0011918951132299

$ws.Range("Q5").Value2 = 1.4022947024663319
$ws.Range("R5").Value2 = 51.528864223543721
$ws.Range("U5").Value2 = 1.4022947024663319
$ws.Range("V5").Value2 = 51.528864223543721

$ws.Range("N6").Value2 = 0.95454545454545459
$ws.Range("R6").Value2 = 44.559542316018756
$ws.Range("S6").Value2 = 3.5882352941176472
$ws.Range("V6").Value2 = 43.949967998029798

$ws.Range("Q7").Value2 = 1.6661332561106461
$ws.Range("R7").Value2 = 63.012801853344513
$ws.Range("U7").Value2 = 1.6661332561106461
$ws.Range("V7").Value2 = 63.012801853344513

$ws.Range("Q8").Value2 = 1.4158531633614351
$ws.Range("R8").Value2 = 42.603670915964123
$ws.Range("U8").Value2 = 1.4158531633614351
$ws.Range("V8").Value2 = 42.603670915964123

$ws.Range("N9").Value2 = 0.95867768595041325
$ws.Range("Q9").Value2 = 1.4313180424938461
$ws.Range("R9").Value2 = 66.609914385233864
$ws.Range("S9").Value2 = 4.1842105263157894
$ws.Range("T9").Value2 = 4.4117647058823533
$ws.Range("V9").Value2 = 65.534657837676789
$ws.Range("Z9").Value2 = 1.0056242969628799

$ws.Range("Q10").Value2 = 1.4490952623589211
$ws.Range("R10").Value2 = 48.874427916309131
$ws.Range("S10").Value2 = 4.2592592592592604
$ws.Range("T10").Value2 = 4.2592592592592604
$ws.Range("U10").Value2 = 1.4490952623589211
$ws.Range("V10").Value2 = 48.874427916309131

$ws.Range("R11").Value2 = 70.030465564385864
$ws.Range("V11").Value2 = 70.030465564385864

# --- New summary row just under the data table: average of the k column ---
$ws.Range("J12").Formula = "=AVERAGE(J2:J11)"
$ws.Range("J12").Font.Bold = $true

# --- New summary block (rows 14-17): average/worst of the two ratio columns ---
$ws.Range("A14").Value = "Average of SW(S*)/SW(OPT)"
$ws.Range("B14").Formula = "=AVERAGE(N2:N11)"

$ws.Range("A15").Value = "Average of SC(S*)/SC(OPT)"
$ws.Range("B15").Formula = "=AVERAGE(Z2:Z11)"

$ws.Range("A16").Value = "Worst of SW(S*)/SW(OPT)"
$ws.Range("B16").Formula = "=MIN(N2:N11)"

$ws.Range("A17").Value = "Worst of SC(S*)/SC(OPT)"
$ws.Range("B17").Formula = "=MAX(Z2:Z11)"

$ws.Range("B14:B17").Font.Bold = $true
$ws.Range("B14:B17").Font.Size = 12
$ws.Range("B14:B17").VerticalAlignment = -4108
